$d = $word.ActiveDocument

# Update the title/date paragraph
$p = $d.Paragraphs.Item(1)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "2025-11-06 Thursday"

# Update each math-problem cell in the table
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "89-32="
$t.Cell(1, 2).Range.Text = "74-73="
$t.Cell(1, 3).Range.Text = "95-37="
$t.Cell(1, 4).Range.Text = "95-25="
$t.Cell(1, 5).Range.Text = "19+71="
$t.Cell(2, 1).Range.Text = "20+26="
$t.Cell(2, 2).Range.Text = "6+20="
$t.Cell(2, 3).Range.Text = "7+32="
$t.Cell(2, 4).Range.Text = "78-19="
$t.Cell(2, 5).Range.Text = "89-15="
$t.Cell(3, 1).Range.Text = "25-7="
$t.Cell(3, 2).Range.Text = "25+53="
$t.Cell(3, 3).Range.Text = "20+24="
$t.Cell(3, 4).Range.Text = "50-15="
$t.Cell(3, 5).Range.Text = "41-31="
$t.Cell(4, 1).Range.Text = "63-2="
$t.Cell(4, 2).Range.Text = "25+41="
$t.Cell(4, 3).Range.Text = "98-30="
$t.Cell(4, 4).Range.Text = "72-25="
$t.Cell(4, 5).Range.Text = "75-10="
$t.Cell(5, 1).Range.Text = "73-13="
$t.Cell(5, 2).Range.Text = "54-8="
$t.Cell(5, 3).Range.Text = "33-0="
$t.Cell(5, 4).Range.Text = "70+10="
$t.Cell(5, 5).Range.Text = "79-28="
$t.Cell(6, 1).Range.Text = "13+21="
$t.Cell(6, 2).Range.Text = "30-17="
$t.Cell(6, 3).Range.Text = "68+22="
$t.Cell(6, 4).Range.Text = "86-3="
$t.Cell(6, 5).Range.Text = "0+61="
$t.Cell(7, 1).Range.Text = "17+50="
$t.Cell(7, 2).Range.Text = "87-45="
$t.Cell(7, 3).Range.Text = "93-8="
$t.Cell(7, 4).Range.Text = "55+9="
$t.Cell(7, 5).Range.Text = "85-21="
$t.Cell(8, 1).Range.Text = "55+27="
$t.Cell(8, 2).Range.Text = "29+36="
$t.Cell(8, 3).Range.Text = "79-72="
$t.Cell(8, 4).Range.Text = "82-63="
$t.Cell(8, 5).Range.Text = "45+8="
$t.Cell(9, 1).Range.Text = "7+81="
$t.Cell(9, 2).Range.Text = "7-4="
$t.Cell(9, 3).Range.Text = "91-34="
$t.Cell(9, 4).Range.Text = "53-43="
$t.Cell(9, 5).Range.Text = "49-23="
$t.Cell(10, 1).Range.Text = "59+22="
$t.Cell(10, 2).Range.Text = "47-13="
$t.Cell(10, 3).Range.Text = "92-23="
$t.Cell(10, 4).Range.Text = "62-38="
$t.Cell(10, 5).Range.Text = "58+27="
$t.Cell(11, 1).Range.Text = "88-70="
$t.Cell(11, 2).Range.Text = "86-1="
$t.Cell(11, 3).Range.Text = "81+5="
$t.Cell(11, 4).Range.Text = "56-41="
$t.Cell(11, 5).Range.Text = "83-65="
$t.Cell(12, 1).Range.Text = "61-6="
$t.Cell(12, 2).Range.Text = "75-44="
$t.Cell(12, 3).Range.Text = "40-16="
$t.Cell(12, 4).Range.Text = "58-5="
$t.Cell(12, 5).Range.Text = "85+11="
$t.Cell(13, 1).Range.Text = "99-55="
$t.Cell(13, 2).Range.Text = "42+20="
$t.Cell(13, 3).Range.Text = "13+49="
$t.Cell(13, 4).Range.Text = "2+66="
$t.Cell(13, 5).Range.Text = "41+11="
$t.Cell(14, 1).Range.Text = "10+74="
$t.Cell(14, 2).Range.Text = "86-20="
$t.Cell(14, 3).Range.Text = "26+11="
$t.Cell(14, 4).Range.Text = "16+12="
$t.Cell(14, 5).Range.Text = "79+10="
$t.Cell(15, 1).Range.Text = "95+3="
$t.Cell(15, 2).Range.Text = "33+22="
$t.Cell(15, 3).Range.Text = "28+28="
$t.Cell(15, 4).Range.Text = "37+32="
$t.Cell(15, 5).Range.Text = "99-55="
$t.Cell(16, 1).Range.Text = "71+0="
$t.Cell(16, 2).Range.Text = "62-32="
$t.Cell(16, 3).Range.Text = "25+45="
$t.Cell(16, 4).Range.Text = "1+47="
$t.Cell(16, 5).Range.Text = "8-4="
$t.Cell(17, 1).Range.Text = "61-61="
$t.Cell(17, 2).Range.Text = "61+27="
$t.Cell(17, 3).Range.Text = "84-14="
$t.Cell(17, 4).Range.Text = "42+34="
$t.Cell(17, 5).Range.Text = "18+27="
$t.Cell(18, 1).Range.Text = "77+11="
$t.Cell(18, 2).Range.Text = "62-3="
$t.Cell(18, 3).Range.Text = "7+48="
$t.Cell(18, 4).Range.Text = "27+47="
$t.Cell(18, 5).Range.Text = "36+52="
$t.Cell(19, 1).Range.Text = "54-42="
$t.Cell(19, 2).Range.Text = "88-65="
$t.Cell(19, 3).Range.Text = "94-34="
$t.Cell(19, 4).Range.Text = "42-19="
$t.Cell(19, 5).Range.Text = "61+31="
$t.Cell(20, 1).Range.Text = "78+17="
$t.Cell(20, 2).Range.Text = "83-18="
$t.Cell(20, 3).Range.Text = "89-9="
$t.Cell(20, 4).Range.Text = "12+69="
$t.Cell(20, 5).Range.Text = "31+24="
